$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are written as exact text (matching
# the source workbook, which stores these as inline strings, not numbers).
$textCells = @(
    "D2",
    "D3",
    "D4",
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D28",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50",
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.49"
$ws.Range("D3").Value = "23.95"
$ws.Range("D4").Value = "5.241"
$ws.Range("D5").Value = "0.05864"
$ws.Range("D7").Value = "3.332"
$ws.Range("D8").Value = "0.8081"
$ws.Range("D9").Value = "0.8833"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01046"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1377"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07251"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03055"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03051"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09323"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.846"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001552"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04710"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.006247"
$ws.Range("D20").Value = "0.001262"
$ws.Range("D21").Value = "0.004582"
$ws.Range("D22").Value = "0.00008698"
$ws.Range("D23").Value = "3.560"
$ws.Range("D24").Value = "2.180"
$ws.Range("D26").Value = "0.1315"
$ws.Range("D28").Value = "0.0002339"
$ws.Range("D41").Value = "0.006378"
$ws.Range("D42").Value = "0.1052"
$ws.Range("D43").Value = "0.002699"
$ws.Range("D44").Value = "0.007811"
$ws.Range("D45").Value = "0.00005481"
$ws.Range("D47").Value = "0.5398"
$ws.Range("D48").Value = "0.002369"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").Value = "0.0001999"

Write-Host "Applied 61 cell updates"
